# CIDC-1278 first pass at redone docs
# Renames "RNAseq Analysis" tab to "RNAseq level 1 Analysis", updates related
# labels/legend text, makes "Legend" the active/selected tab, and normalizes
# the selections on the other sheets.

$wb = $excel.ActiveWorkbook

$wsRnaseq   = $wb.Worksheets.Item(1)   # "RNAseq Analysis"
$wsExcluded = $wb.Worksheets.Item(2)   # "Excluded Samples"
$wsLegend   = $wb.Worksheets.Item(3)   # "Legend"
$wsDataDict = $wb.Worksheets.Item(4)   # "Data Dictionary"

# --- Rename the first sheet ---
$wsRnaseq.Name = "RNAseq level 1 Analysis"

# --- Update text content that references the old section/tab names ---
$wsRnaseq.Range("B5").Value = "RNAseq level 1 Runs"
$wsLegend.Range("B5").Value = "Section 'RNAseq level 1 Runs' of tab 'RNAseq level 1Analysis'"

# --- Normalize selection on "RNAseq level 1 Analysis" (no longer the active tab) ---
$wsRnaseq.Range("B6").Select()

# --- Normalize selection on "Excluded Samples" ---
$wsExcluded.Range("B3").Select()

# --- Normalize selection on "Data Dictionary" ---
$wsDataDict.Range("A1").Select()

# --- Make "Legend" the active tab with a single-cell selection ---
$wsLegend.Activate()
$wsLegend.Range("B6").Select()

# --- Adjust auto row heights on the "Legend" sheet (wrapped text rows) ---
$wsLegend.Rows.Item(2).RowHeight = 13.8
$wsLegend.Rows.Item(5).RowHeight = 23.95
